# Apply crypto price/volume update for Mon Jan 30 21:52:11 UTC 2023 GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B/C (plain text coin name + link) swap for rows 7 and 8 -- ordinary text assignment is fine here
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("B8").Value = "FTXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"

# D/E (Price / Volume) columns hold numeric-looking text ("305.72", "-3.82%") that must stay
# plain text (matches the source inlineStr cells), not get auto-converted to Number/Percent by
# Excel's normal type inference. Force text format on the full D2:E51 data block first, assign
# every changed value as a string, then clear the temporary formatting so the cells end up back
# at the default (unstyled) General format -- exactly like the untouched cells around them.
$priceVolRange = $ws.Range("D2:E51")
$priceVolRange.NumberFormat = "@"

$ws.Range("D2").Value = "305.72"
$ws.Range("E2").Value = "-3.82%"
$ws.Range("D3").Value = "37.24"
$ws.Range("E3").Value = "-6.28%"
$ws.Range("D4").Value = "5.089"
$ws.Range("E4").Value = "-0.92%"
$ws.Range("D5").Value = "0.07712"
$ws.Range("E5").Value = "-6.04%"
$ws.Range("D6").Value = "4.355"
$ws.Range("E6").Value = "1.24%"
$ws.Range("D7").Value = "8.199"
$ws.Range("E7").Value = "-1.92%"
$ws.Range("D8").Value = "1.879"
$ws.Range("E8").Value = "-8.82%"
$ws.Range("D9").Value = "3.125"
$ws.Range("E9").Value = "-3.44%"
$ws.Range("D10").Value = "0.9172"
$ws.Range("E10").Value = "-2.17%"
$ws.Range("D11").Value = "0.1160"
$ws.Range("E11").Value = "-14.41%"
$ws.Range("D12").Value = "0.1875"
$ws.Range("E12").Value = "-5.19%"
$ws.Range("D13").Value = "0.08688"
$ws.Range("E13").Value = "-4.18%"
$ws.Range("D14").Value = "0.03411"
$ws.Range("E14").Value = "-2.62%"
$ws.Range("D15").Value = "0.09696"
$ws.Range("E15").Value = "-0.94%"
$ws.Range("D16").Value = "0.001369"
$ws.Range("E16").Value = "-2.51%"
$ws.Range("D17").Value = "0.005923"
$ws.Range("E17").Value = "-4.75%"
$ws.Range("D18").Value = "3.587"
$ws.Range("E18").Value = "-2.56%"
$ws.Range("D19").Value = "0.3407"
$ws.Range("E19").Value = "-2.08%"
$ws.Range("D20").Value = "0.1276"
$ws.Range("E20").Value = "-3.63%"
$ws.Range("D21").Value = "5.023"
$ws.Range("E21").Value = "0.65%"
$ws.Range("D23").Value = "0.02108"
$ws.Range("E23").Value = "5,173.45%"
$ws.Range("D24").Value = "0.04331"
$ws.Range("E24").Value = "-0.43%"
$ws.Range("D25").Value = "0.001215"
$ws.Range("E25").Value = "-1.12%"
$ws.Range("D26").Value = "0.004535"
$ws.Range("D27").Value = "0.0001354"
$ws.Range("E27").Value = "3.99%"
$ws.Range("D39").Value = "0.02192"
$ws.Range("E39").Value = "-3.12%"
$ws.Range("D40").Value = "0.04897"
$ws.Range("E40").Value = "-5.77%"
$ws.Range("D41").Value = "0.007555"
$ws.Range("E41").Value = "-2.68%"
$ws.Range("D42").Value = "0.009909"
$ws.Range("E42").Value = "0.60%"
$ws.Range("D43").Value = "0.1334"
$ws.Range("E43").Value = "-5.17%"
$ws.Range("D44").Value = "0.002066"
$ws.Range("E44").Value = "0.98%"
$ws.Range("D45").Value = "0.008515"
$ws.Range("E45").Value = "-12.00%"
$ws.Range("E46").Value = "-0.81%"
$ws.Range("D47").Value = "0.00000000752"
$ws.Range("E47").Value = "0.28%"
$ws.Range("D48").Value = "0.003002"
$ws.Range("E48").Value = "1.92%"
$ws.Range("D49").Value = "0.001303"
$ws.Range("E49").Value = "-22.83%"
$ws.Range("D50").Value = "0.00002105"
$ws.Range("E50").Value = "0.28%"
$ws.Range("D51").Value = "0.0002005"
$ws.Range("E51").Value = "0.28%"

$priceVolRange.ClearFormats()
